# Add a new column C ("activities count" per user row): 2 for the first
# 8 users, 4 for the remaining 7 (CRUD activity counters per user row).
$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1:C8").Value = 2
$ws.Range("C9:C15").Value = 4

# Leave the whole populated table selected, matching the post-edit view state.
$ws.Range("A1:C15").Select()
